# Apply the latest cryptos.xlsx price/volume/ranking refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell as TEXT (not auto-converted to a number),
# matching the source data which stores prices like "1.000" / "29.988.46"
# as plain strings rather than numerics.
function Set-TextValue($cellRef, $text) {
    $ws.Range($cellRef).Value = "'" + $text
}

Set-TextValue 'D2' '29.988.46'
$ws.Range('E2').Value = '  -0.21%  '

Set-TextValue 'D3' '1.897.60'

Set-TextValue 'D4' '1.000'
$ws.Range('E4').Value = '  +0.10%  '

Set-TextValue 'D5' '0.8411'
$ws.Range('E5').Value = '  +3.17%  '

Set-TextValue 'D6' '241.72'
$ws.Range('E6').Value = '  -0.62%  '

Set-TextValue 'D7' '1.001'
$ws.Range('E7').Value = '  +0.08%  '

Set-TextValue 'D8' '0.3313'
$ws.Range('E8').Value = '  +3.48%  '

Set-TextValue 'D9' '26.66'
$ws.Range('E9').Value = '  +0.78%  '

Set-TextValue 'D10' '0.07057'
$ws.Range('E10').Value = '  +1.38%  '

Set-TextValue 'D11' '0.08079'
$ws.Range('E11').Value = '  +0.49%  '

Set-TextValue 'D12' '0.7610'
$ws.Range('E12').Value = '  +0.95%  '

Set-TextValue 'D13' '1.898.88'
$ws.Range('E13').Value = '  -0.68%  '

Set-TextValue 'D14' '5.266'
$ws.Range('E14').Value = '  +0.35%  '

Set-TextValue 'D15' '92.27'
$ws.Range('E15').Value = '  -1.76%  '

Set-TextValue 'D16' '29.991.72'
$ws.Range('E16').Value = '  -0.21%  '

$ws.Range('E17').Value = '  +0.26%  '

Set-TextValue 'D18' '5.882'
$ws.Range('E18').Value = '  -2.54%  '

Set-TextValue 'D19' '244.41'
$ws.Range('E19').Value = '  -2.65%  '

Set-TextValue 'D20' '0.000007769'
$ws.Range('E20').Value = '  -0.50%  '

Set-TextValue 'D21' '0.9997'
$ws.Range('E21').Value = '  +0.01%  '

Set-TextValue 'D22' '2.149.08'
$ws.Range('E22').Value = '  -0.36%  '

Set-TextValue 'D23' '1.000'
$ws.Range('E23').Value = '  +0.12%  '

Set-TextValue 'D24' '6.980'
$ws.Range('E24').Value = '  -0.25%  '

Set-TextValue 'D25' '0.1752'
$ws.Range('E25').Value = '  +23.27%  '

Set-TextValue 'D26' '9.259'
$ws.Range('E26').Value = '  -0.97%  '

Set-TextValue 'D27' '166.31'
$ws.Range('E27').Value = '  -1.47%  '

Set-TextValue 'D28' '18.90'
$ws.Range('E28').Value = '  -0.84%  '

Set-TextValue 'D29' '2.106'
$ws.Range('E29').Value = '  +1.80%  '

$ws.Range('E30').Value = '  -2.30%  '

Set-TextValue 'D31' '1.520'
$ws.Range('E31').Value = '  -0.42%  '

Set-TextValue 'D32' '0.05862'
$ws.Range('E32').Value = '  +9.09%  '

$ws.Range('E33').Value = '  -2.01%  '

Set-TextValue 'D34' '4.080'
$ws.Range('E34').Value = '  -1.24%  '

Set-TextValue 'D35' '1.272'
$ws.Range('E35').Value = '  +0.24%  '

Set-TextValue 'D36' '0.7314'
$ws.Range('E36').Value = '  -1.38%  '

Set-TextValue 'D37' '2.720'
$ws.Range('E37').Value = '  -0.22%  '

Set-TextValue 'D38' '0.01919'
$ws.Range('E38').Value = '  -0.92%  '

Set-TextValue 'D39' '2.773'
$ws.Range('E39').Value = '  -0.44%  '

$ws.Range('E40').Value = '  -1.09%  '

Set-TextValue 'D41' '72.56'
$ws.Range('E41').Value = '  -0.73%  '

Set-TextValue 'D42' '5.853'
$ws.Range('E42').Value = '  -5.45%  '

Set-TextValue 'D43' '0.8445'
$ws.Range('E43').Value = '  +1.22%  '

Set-TextValue 'D44' '1.000'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('E45').Value = '  -1.44%  '

Set-TextValue 'D46' '101.67'
$ws.Range('E46').Value = '  +0.71%  '

Set-TextValue 'D47' '1.014.80'
$ws.Range('E47').Value = '  +5.36%  '

$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '9.853'
$ws.Range('E48').Value = '  -0.09%  '

$ws.Range('B49').Value = 'Aptos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D49' '7.572'
$ws.Range('E49').Value = '  -1.41%  '

Set-TextValue 'D50' '2.047.67'
$ws.Range('E50').Value = '  -0.41%  '

Set-TextValue 'D51' '35.92'
$ws.Range('E51').Value = '  -1.62%  '
